# [AFG] added final excel sheets for Afghanistan
#
# 1. Remove the 4 stray empty cells in "ODI Batting" (B3, B5, B6, B7).
# 2. Add two new worksheets "ODI Batting Extra" and "ODI Bowling Extra"
#    (after "ODI Bowling"), each with a bold/centred header row copied
#    from the existing header style, and their data rows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) ODI Batting: drop the leftover empty inline-string cells in column B
# ---------------------------------------------------------------------
$odiBatting = $wb.Worksheets.Item("ODI Batting")
$odiBatting.Range("B3").ClearContents()
$odiBatting.Range("B5").ClearContents()
$odiBatting.Range("B6").ClearContents()
$odiBatting.Range("B7").ClearContents()

# ---------------------------------------------------------------------
# 2) Add "ODI Batting Extra" right after "ODI Bowling"
# ---------------------------------------------------------------------
$odiBowling = $wb.Worksheets.Item("ODI Bowling")
$battingExtra = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $odiBowling)
$battingExtra.Name = "ODI Batting Extra"

$battingExtraHeaders = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($i = 0; $i -lt $battingExtraHeaders.Length; $i++) {
    $battingExtra.Cells.Item(1, $i + 1).Value = $battingExtraHeaders[$i]
}
# Re-use the existing bold/centred header style (style index 1) instead of
# fabricating a new one.
$odiBatting.Range("A1").Copy()
$battingExtra.Range("A1:F1").PasteSpecial(-4122)

# MATCH_CODE, BATTING_POSITION, NUM_4, NUM_6, PERCENT_RUNS_OF_TOTAL, MAN_OF_MATCH
$battingExtraRows = @(
    @("3864", 11, "0", "0", $null, "NO"),
    @("3917", 11, $null, $null, $null, "NO"),
    @("3918", 11, $null, $null, $null, "NO"),
    @("4379", $null, $null, $null, $null, "NO"),
    @("4525", 11, $null, $null, $null, "NO"),
    @("4528", 10, "0", "0", "2.33%", "NO"),
    @("4537", $null, $null, $null, $null, "NO"),
    @("4671", 11, "0", "0", "0.44%", "NO"),
    @("4674", $null, $null, $null, $null, $null)
)

for ($r = 0; $r -lt $battingExtraRows.Length; $r++) {
    $row = $battingExtraRows[$r]
    $xlRow = $r + 2

    # MATCH_CODE (A) - numeric-looking but stored as text
    $battingExtra.Cells.Item($xlRow, 1).Value = "'" + $row[0]

    # BATTING_POSITION (B) - real number
    if ($null -ne $row[1]) {
        $battingExtra.Cells.Item($xlRow, 2).Value = $row[1]
    }

    # NUM_4 (C) / NUM_6 (D) - numeric-looking text
    if ($null -ne $row[2]) {
        $battingExtra.Cells.Item($xlRow, 3).Value = "'" + $row[2]
    }
    if ($null -ne $row[3]) {
        $battingExtra.Cells.Item($xlRow, 4).Value = "'" + $row[3]
    }

    # PERCENT_RUNS_OF_TOTAL (E) - percentage-looking text
    if ($null -ne $row[4]) {
        $battingExtra.Cells.Item($xlRow, 5).Value = "'" + $row[4]
    }

    # MAN_OF_MATCH (F) - plain text, no trick needed
    if ($null -ne $row[5]) {
        $battingExtra.Cells.Item($xlRow, 6).Value = $row[5]
    }
}

# ---------------------------------------------------------------------
# 3) Add "ODI Bowling Extra" right after "ODI Batting Extra"
# ---------------------------------------------------------------------
$bowlingExtra = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $battingExtra)
$bowlingExtra.Name = "ODI Bowling Extra"

$bowlingExtraHeaders = @("MATCH_CODE", "MAIDEN_OVERS", "PERCENT_WICKETS_OF_ALL")
for ($i = 0; $i -lt $bowlingExtraHeaders.Length; $i++) {
    $bowlingExtra.Cells.Item(1, $i + 1).Value = $bowlingExtraHeaders[$i]
}
$odiBatting.Range("A1").Copy()
$bowlingExtra.Range("A1:C1").PasteSpecial(-4122)

# MATCH_CODE, MAIDEN_OVERS, PERCENT_WICKETS_OF_ALL
$bowlingExtraRows = @(
    @("3864", "1", $null),
    @("3917", "1", "10.00%"),
    @("3918", "0", "10.00%"),
    @("4379", "1", $null),
    @("4525", $null, $null),
    @("4528", "0", "10.00%"),
    @("4537", "0", $null),
    @("4671", $null, $null)
)

for ($r = 0; $r -lt $bowlingExtraRows.Length; $r++) {
    $row = $bowlingExtraRows[$r]
    $xlRow = $r + 2

    # MATCH_CODE (A) - numeric-looking but stored as text
    $bowlingExtra.Cells.Item($xlRow, 1).Value = "'" + $row[0]

    # MAIDEN_OVERS (B) - numeric-looking text
    if ($null -ne $row[1]) {
        $bowlingExtra.Cells.Item($xlRow, 2).Value = "'" + $row[1]
    }

    # PERCENT_WICKETS_OF_ALL (C) - percentage-looking text
    if ($null -ne $row[2]) {
        $bowlingExtra.Cells.Item($xlRow, 3).Value = "'" + $row[2]
    }
}
